# Auto-update draw results: append the 2025-11-23 Pick 3 result as a new
# row at the bottom of the "Results" sheet (row 68), mirroring the layout
# of every prior row: Date, Game, Phase, Result, InsertedAt - all stored
# as literal text (not auto-coerced into dates/numbers by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 68

$date       = "2025-11-23"
$game       = "Pick 3"
$phase      = "251123"
$result     = "5-6-1"
$insertedAt = "2025-11-23T21:37:33.156+04:00"

# Column A ("Date") and C ("Phase") look like a date / a plain number to
# Excel's auto-detection, so force the cell to Text format first (then
# restore the default "Normal" style after writing so no stray number
# format lingers on the cell) - this keeps the value as the literal
# string instead of being reinterpreted as a date serial or a number.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = $date
$cellA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = $game

$cellC = $ws.Cells.Item($newRow, 3)
$cellC.NumberFormat = "@"
$cellC.Value = $phase
$cellC.Style = "Normal"

$ws.Cells.Item($newRow, 4).Value = $result
$ws.Cells.Item($newRow, 5).Value = $insertedAt

# Extend the "number stored as text" ignored-error range to cover the
# freshly appended row, matching the rest of the sheet's A1:E<lastRow>
# ignore range.
$ws.Range("A1:E" + $newRow).Errors.Item(9).Ignore = $true
